$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free direct cell updates. For values that look numeric (e.g. "109.70"),
# force text storage via NumberFormat "@" so Excel does not coerce them to doubles
# (which would silently drop formatting like trailing zeros or thousands separators).

$ws.Range('D2').Value = '51.954.42'
$ws.Range('E2').Value = '  -0.57%  '
$ws.Range('D3').Value = '2.791.10'
$ws.Range('E3').Value = '  -1.98%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '360.24'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.70'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.46%  '
$ws.Range('E7').Value = '  -2.70%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -2.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.17'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0849'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.55%  '
$ws.Range('E12').Value = '  +1.22%  '
$ws.Range('E13').Value = '  -2.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.57'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.84%  '
$ws.Range('D15').Value = '3.227.87'
$ws.Range('E15').Value = '  -2.04%  '
$ws.Range('D16').Value = '2.809.09'
$ws.Range('E16').Value = '  -1.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.940'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +3.93%  '
$ws.Range('D18').Value = '51.915.33'
$ws.Range('E18').Value = '  -0.52%  '
$ws.Range('E19').Value = '  -1.59%  '
$ws.Range('E20').Value = '  -2.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.13'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.60%  '
$ws.Range('E22').Value = '  -1.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.43'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '269.77'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.76'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.56'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.52%  '
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.160'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +13.89%  '
$ws.Range('E29').Value = '  -1.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.29'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0471'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.76%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '51.98'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.44'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('E34').Value = '  -2.77%  '
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.25'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -3.12%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.03'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +3.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.20'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.19%  '
$ws.Range('E40').Value = '  -4.00%  '
$ws.Range('E41').Value = '  +1.96%  '
$ws.Range('E42').Value = '  -2.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.24'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.57%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.03'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -7.77%  '
$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '119.47'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -7.10%  '
$ws.Range('D46').Value = '2.084.05'
$ws.Range('E46').Value = '  -1.67%  '
$ws.Range('E47').Value = '  -4.63%  '
$ws.Range('E49').Value = '  -1.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.956'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -5.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.84'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.67%  '
